# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to specific Leve rows across ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5443.6
$ws.Range("I32").Value = 6249.5
$ws.Range("J32").Value = 4906.3335
$ws.Range("K32").Value = 6249.5
$ws.Range("L32").Value = 4906.3335
$ws.Range("M32").Value = -5923.5
$ws.Range("N32").Value = -5558.3335

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H74").Value = 2247.5
$ws.Range("I74").Value = 995
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 995
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -59

$ws.Range("H77").Value = 2247.5
$ws.Range("I77").Value = 995
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 4975
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -295

$ws.Range("H98").Value = 1660.3334
$ws.Range("I98").Value = 1827.25
$ws.Range("J98").Value = 1326.5
$ws.Range("K98").Value = 1827.25
$ws.Range("L98").Value = 1326.5
$ws.Range("M98").Value = -329.25
$ws.Range("N98").Value = -4322.5

$ws.Range("H103").Value = 1539.8
$ws.Range("I103").Value = 2349.5
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 7048.5
$ws.Range("L103").Value = 3000
$ws.Range("M103").Value = -6462.5
$ws.Range("N103").Value = -4172

$ws.Range("H106").Value = 6990
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 6990
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 6990
$ws.Range("N106").Value = -8252

$ws.Range("H122").Value = 1660.3334
$ws.Range("I122").Value = 1827.25
$ws.Range("J122").Value = 1326.5
$ws.Range("K122").Value = 5481.75
$ws.Range("L122").Value = 3979.5
$ws.Range("M122").Value = -3031.75
$ws.Range("N122").Value = -8879.5

$ws.Range("H132").Value = 939.1053000000001
$ws.Range("I132").Value = 1019.0714
$ws.Range("J132").Value = 715.2
$ws.Range("K132").Value = 3057.2142
$ws.Range("L132").Value = 2145.6
$ws.Range("M132").Value = -527.2142000000003

$ws.Range("H141").Value = 8665
$ws.Range("I141").Value = 7997.5
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 23992.5
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -18812.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5998
$ws.Range("I61").Value = 6331.1665
$ws.Range("J61").Value = 3999
$ws.Range("K61").Value = 6331.1665
$ws.Range("L61").Value = 3999
$ws.Range("M61").Value = -6119.1665
$ws.Range("N61").Value = -4423

$ws.Range("H88").Value = 5000
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -4594
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 5000
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -3596
$ws.Range("N91").ClearContents()

$ws.Range("H102").Value = 776.8
$ws.Range("I102").Value = 776.8
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 776.8
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 845.2

$ws.Range("H136").Value = 5998
$ws.Range("I136").Value = 6331.1665
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 18993.4995
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -16443.4995
$ws.Range("N136").Value = -17097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3136.75
$ws.Range("I20").Value = 3357
$ws.Range("J20").Value = 2769.6667
$ws.Range("K20").Value = 3357
$ws.Range("L20").Value = 2769.6667
$ws.Range("M20").Value = -3110
$ws.Range("N20").Value = -3263.6667

$ws.Range("H64").Value = 1199.2
$ws.Range("I64").Value = 999
$ws.Range("J64").Value = 1249.25
$ws.Range("K64").Value = 999
$ws.Range("L64").Value = 1249.25
$ws.Range("M64").Value = -774
$ws.Range("N64").Value = -1699.25

$ws.Range("H67").Value = 1199.2
$ws.Range("I67").Value = 999
$ws.Range("J67").Value = 1249.25
$ws.Range("K67").Value = 999
$ws.Range("L67").Value = 1249.25
$ws.Range("M67").Value = -219
$ws.Range("N67").Value = -2809.25

$ws.Range("H86").Value = 2819.7144
$ws.Range("I86").Value = 3047.6
$ws.Range("J86").Value = 2250
$ws.Range("K86").Value = 3047.6
$ws.Range("L86").Value = 2250
$ws.Range("M86").Value = -1924.6
$ws.Range("N86").Value = -4496

$ws.Range("H89").Value = 2819.7144
$ws.Range("I89").Value = 3047.6
$ws.Range("J89").Value = 2250
$ws.Range("K89").Value = 15238
$ws.Range("L89").Value = 11250
$ws.Range("M89").Value = -9622
$ws.Range("N89").Value = -22482

$ws.Range("H94").Value = 2132.8333
$ws.Range("I94").Value = 2474.5
$ws.Range("J94").Value = 1449.5
$ws.Range("K94").Value = 2474.5
$ws.Range("L94").Value = 1449.5
$ws.Range("M94").Value = -2023.5
$ws.Range("N94").Value = -2351.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4261
$ws.Range("I31").Value = 1265.4615
$ws.Range("J31").Value = 13996.5
$ws.Range("K31").Value = 1265.4615
$ws.Range("L31").Value = 13996.5
$ws.Range("M31").Value = -970.4614999999999

$ws.Range("H34").Value = 4261
$ws.Range("I34").Value = 1265.4615
$ws.Range("J34").Value = 13996.5
$ws.Range("K34").Value = 1265.4615
$ws.Range("L34").Value = 13996.5
$ws.Range("M34").Value = -1063.4615

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 662.6667
$ws.Range("I12").Value = 719.5
$ws.Range("J12").Value = 549
$ws.Range("K12").Value = 2158.5
$ws.Range("L12").Value = 1647
$ws.Range("M12").Value = -1985.5
$ws.Range("N12").Value = -1993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H102").Value = 4516.3335
$ws.Range("I102").Value = 4516.3335
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4516.3335
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2894.3335

$ws.Range("H107").Value = 625
$ws.Range("I107").Value = 720
$ws.Range("J107").Value = 435
$ws.Range("K107").Value = 720
$ws.Range("L107").Value = 435
$ws.Range("M107").Value = 1200
$ws.Range("N107").Value = -4275

$ws.Range("H113").Value = 2412.8333
$ws.Range("I113").Value = 2326
$ws.Range("J113").Value = 2499.6667
$ws.Range("K113").Value = 2326
$ws.Range("L113").Value = 2499.6667
$ws.Range("M113").Value = -156
$ws.Range("N113").Value = -6839.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4145.75
$ws.Range("I7").Value = 4145.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4145.75
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4033.75

$ws.Range("H126").Value = 4145.75
$ws.Range("I126").Value = 4145.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12437.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9967.25
